# Update "想去人数" (F column) figures across sheets, per refreshed data pull.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(5, 6).Value = 4515
$ws1.Cells.Item(6, 6).Value = 1810
$ws1.Cells.Item(7, 6).Value = 121
$ws1.Cells.Item(9, 6).Value = 3051
$ws1.Cells.Item(10, 6).Value = 581
$ws1.Cells.Item(11, 6).Value = 579
$ws1.Cells.Item(13, 6).Value = 581
$ws1.Cells.Item(14, 6).Value = 502
$ws1.Cells.Item(15, 6).Value = 500
$ws1.Cells.Item(17, 6).Value = 126
$ws1.Cells.Item(18, 6).Value = 1748
$ws1.Cells.Item(19, 6).Value = 1286
$ws1.Cells.Item(21, 6).Value = 1521
$ws1.Cells.Item(22, 6).Value = 121
$ws1.Cells.Item(23, 6).Value = 597
$ws1.Cells.Item(24, 6).Value = 39
$ws1.Cells.Item(25, 6).Value = 523
$ws1.Cells.Item(27, 6).Value = 35
$ws1.Cells.Item(28, 6).Value = 82
$ws1.Cells.Item(29, 6).Value = 116
$ws1.Cells.Item(31, 6).Value = 3376
$ws1.Cells.Item(33, 6).Value = 58
$ws1.Cells.Item(34, 6).Value = 220
$ws1.Cells.Item(35, 6).Value = 51
$ws1.Cells.Item(36, 6).Value = 1648

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 19
$ws2.Cells.Item(3, 6).Value = 32

# --- Sheet "全部类型" (all types, union of the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 4515
$ws4.Cells.Item(6, 6).Value = 1810
$ws4.Cells.Item(7, 6).Value = 121
$ws4.Cells.Item(9, 6).Value = 3051
$ws4.Cells.Item(10, 6).Value = 581
$ws4.Cells.Item(11, 6).Value = 579
$ws4.Cells.Item(13, 6).Value = 581
$ws4.Cells.Item(14, 6).Value = 502
$ws4.Cells.Item(15, 6).Value = 500
$ws4.Cells.Item(16, 6).Value = 19
$ws4.Cells.Item(18, 6).Value = 126
$ws4.Cells.Item(19, 6).Value = 1748
$ws4.Cells.Item(20, 6).Value = 1286
$ws4.Cells.Item(22, 6).Value = 1521
$ws4.Cells.Item(23, 6).Value = 121
$ws4.Cells.Item(24, 6).Value = 597
$ws4.Cells.Item(25, 6).Value = 39
$ws4.Cells.Item(26, 6).Value = 523
$ws4.Cells.Item(28, 6).Value = 35
$ws4.Cells.Item(29, 6).Value = 82
$ws4.Cells.Item(30, 6).Value = 116
$ws4.Cells.Item(32, 6).Value = 3376
$ws4.Cells.Item(33, 6).Value = 32
$ws4.Cells.Item(35, 6).Value = 58
$ws4.Cells.Item(36, 6).Value = 220
$ws4.Cells.Item(37, 6).Value = 51
$ws4.Cells.Item(38, 6).Value = 1648
